# Update Metadata sheet: Title and Date values
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B5").Value = "NG-Imm Next Dose Date"
$meta.Range("B8").Value = "2025-06-24T09:13:37+01:00"

# Update Elements sheet
$els = $wb.Worksheets.Item("Elements")

# Row 6 (Extension.value[x]) content changes
$els.Range("K6").Value = "date`n"
$els.Range("L6").Value = "Recommended date for the next immunization dose"
$els.Range("AB6").Value = ""
$els.Range("AC6").Value = ""
$els.Range("AE6").Value = ""

# Delete row 7 entirely (Extension.value[x]:valueDate slice, now merged into row 6)
$els.Rows.Item(7).EntireRow.Delete()

# Unhide the data rows
$els.Rows.Item(2).Hidden = $false
$els.Rows.Item(3).Hidden = $false
$els.Rows.Item(4).Hidden = $false
$els.Rows.Item(5).Hidden = $false
$els.Rows.Item(6).Hidden = $false

# Remove the autofilter (data is no longer filtered)
$els.AutoFilterMode = $false

# Clear conditional formatting rules tied to the old filtered view
$els.Cells.FormatConditions.Delete()

# Resize columns to fit the new content
$els.Columns.Item(1).EntireColumn.AutoFit()
$els.Columns.Item(11).EntireColumn.AutoFit()

Write-Output "edit applied"
